$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40-52 down to 41-53.
$ws.Rows(40).Insert()

# Populate the newly inserted row 40 with the new record.
$ws.Cells.Item(40, 1).Value = 5
$ws.Cells.Item(40, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(40, 3).Value = "Maule"
$ws.Cells.Item(40, 4).Value = 44964
$ws.Cells.Item(40, 5).Value = 7
$ws.Cells.Item(40, 6).Value = 100112043
$ws.Cells.Item(40, 7).Value = "Pepino dulce"
$ws.Cells.Item(40, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 500
$ws.Cells.Item(40, 11).Value = 17000
$ws.Cells.Item(40, 12).Value = 17000
$ws.Cells.Item(40, 13).Value = 17000
$ws.Cells.Item(40, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(40, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(40, 16).Value = 944
$ws.Cells.Item(40, 17).Value = 18
$ws.Cells.Item(40, 18).Value = "Hortaliza"
